# Edit script: applies the "Agg. dizionario + messaggio RisultatoCruciverba" change.
#
# 1) Removes the stray _GoBack bookmark that sat between "...se presente" and ".".
# 2) Appends, after the final "1)" paragraph:
#      - a paragraph containing only a page break
#      - a "Risultati ottenuti:" paragraph (Paragrafoelenco style)
#      - a long comment paragraph (Paragrafoelenco + bullet numbering) ending
#        with a fresh _GoBack bookmark (Word always re-anchors _GoBack at the
#        location of the most recent edit).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 1: remove the old _GoBack bookmark without disturbing the surrounding
# run structure. We rebuild the whole paragraph's XML identically, just
# without the bookmarkStart/bookmarkEnd pair.
# ---------------------------------------------------------------------------

$ns = "xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`" xmlns:w14=`"http://schemas.microsoft.com/office/word/2010/wordml`""

$found = $d.Content.Find.Execute("Se è stata trovata una parola")
if (-not $found) {
    throw "Could not locate target paragraph for bookmark removal"
}
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Se è stata trovata una parola*presente.*") {
        $targetPara = $p
    }
}
if ($null -eq $targetPara) {
    throw "Could not find paragraph containing the _GoBack bookmark"
}

$bookmarkParaXml = '<w:p ' + $ns + ' w14:paraId="758F502E" w14:textId="3E540527" w:rsidR="00C81BF6" w:rsidRDefault="00C81BF6" w:rsidP="00C81BF6"><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Se è stata trovata una parola la inserisco nello schema e torno al punto </w:t></w:r><w:r w:rsidR="00884AA7"><w:t>4</w:t></w:r><w:r><w:t xml:space="preserve">), altrimenti torno al punto </w:t></w:r><w:r w:rsidR="006E1B5D"><w:t>4</w:t></w:r><w:r><w:t>)</w:t></w:r><w:r w:rsidR="006E1B5D"><w:t>, prendendo la prossima parola di c caselle se presente</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>'

$targetPara.Range.InsertXML($bookmarkParaXml)

# ---------------------------------------------------------------------------
# Step 2: append the new content after the final "1)" paragraph.
# ---------------------------------------------------------------------------

# The Paragraphs collection reports a handful of trailing zero-length
# "phantom" paragraphs past the real end of body content, so find the last
# paragraph that actually holds text instead of trusting .Last / .Count.
$realLastIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.Trim().Length -gt 0) {
        $realLastIndex = $i
    }
}

$lastPara = $d.Paragraphs.Item($realLastIndex)
$insertion = $lastPara.Range.Duplicate
$insertion.Collapse(0)

# Build three placeholder paragraphs (split via literal paragraph marks) so
# each one becomes a genuine, independently addressable w:p in the OM - then
# fix each one up with the exact target XML.
$insertion.InsertAfter("`rPAGEBREAK`rRisultati ottenuti:`rBODYPLACEHOLDER")

$pBreak = $d.Paragraphs.Item($realLastIndex + 1)
$pTitle = $d.Paragraphs.Item($realLastIndex + 2)
$pBody = $d.Paragraphs.Item($realLastIndex + 3)

$pBreak.Range.InsertXML('<w:p ' + $ns + '><w:r><w:br w:type="page"/></w:r></w:p>')

$pTitle.Range.InsertXML('<w:p ' + $ns + '><w:pPr><w:pStyle w:val="Paragrafoelenco"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Risultati ottenuti:</w:t></w:r></w:p>')

$bodyText = "Credo che valga per tutti gli algoritmi, il dizionario delle parole disponibili/da inserire alla fine dell’elaborazione non è vuoto ma risulta contenere alcune parole le cui lettere sono state tutte inserite durante l’inserimento di altre parole, quindi la parola è completa ma non è stata inserita nello schema tramite la procedura degli algoritmi. Questo fa sì che ci sia una procedura di aggiornamento del dizionario delle parole disponibili che cicla tutte le parole inserite nello schema che sono complete e le elimini dal dizionario se sono contenute al suo interno. AUMENTA LA COMPLESSITà dell’algoritmo/algoritmi."

$pBody.Range.InsertXML('<w:p ' + $ns + '><w:pPr><w:pStyle w:val="Paragrafoelenco"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>' + $bodyText + '</w:t></w:r><w:bookmarkStart w:id="1" w:name="_GoBack"/><w:bookmarkEnd w:id="1"/></w:p>')

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
